# The product name used to contain a space ("...TR-1-Late Repayment").
# It has been renamed to drop the space ("...TR-1-LateRepayment") on both
# the input and the output sheet of the workbook.

$wb = $excel.ActiveWorkbook

$wsInput  = $wb.Worksheets.Item(1)   # ProductLoanInput
$wsOutput = $wb.Worksheets.Item(2)   # ProductLoanOutput

$newProductName = "966-MS-EI-DB-DL-REC-NON-RNI-CTPD-DL-MD-TR-1-LateRepayment"

# Update the product name value on both sheets (cell B1).
$wsInput.Range("B1").Value  = $newProductName
$wsOutput.Range("B1").Value = $newProductName

# Restore/refresh the selection on each sheet without disturbing which
# sheet is the active tab (ProductLoanOutput is the active tab in this
# workbook). Selecting B1 on the input sheet first, then re-activating
# the output sheet and selecting B1 there, mirrors the final cursor
# position recorded for each sheet.
$wsInput.Range("B1").Select()

$wsOutput.Activate()
$wsOutput.Range("B1").Select()
